$d = $word.ActiveDocument

# Locate "150 feet." within the recommendations paragraph and place an
# insertion-point range right after it, then type the new sentence there.
# This mirrors what a human editor does in Word: click right before " It is
# worth noting" and type " Our model suggests..." -- which splits the
# existing run into three runs sharing the same character formatting.

$range = $d.Content
$find = $range.Find
$find.ClearFormatting()
$found = $find.Execute("vertical drop by 150 feet.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insertionPoint = $d.Range($range.End, $range.End)
    $insertionPoint.InsertAfter(" Our model suggests that making these changes supports a ticket price increase of `$1.99.")
}
